$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the replacement data block for rows 2-13 (columns A-T)
$arr = New-Object 'object[,]' 12,20

# Row 2
$arr[0,0] = "ECs"
$arr[0,1] = "Sema3a"
$arr[0,2] = "Plxna2"
$arr[0,3] = "ECs"
$arr[0,4] = [double]3
$arr[0,5] = [double]1
$arr[0,6] = [double]1.339548666666667
$arr[0,7] = [double]4.018646
$arr[0,8] = [double]0.09827984122213275
$arr[0,9] = [double]0.09827984122213274
$arr[0,10] = [double]3
$arr[0,11] = [double]1
$arr[0,12] = [double]39.02618766666667
$arr[0,13] = [double]117.078563
$arr[0,14] = [double]0.4958819606525626
$arr[0,15] = [double]0.4958819606525626
$arr[0,16] = [double]52.27747765396644
$arr[0,17] = [double]470.4972988856981
$arr[0,18] = [double]0.04873520035785373
$arr[0,19] = [double]0.04873520035785373

# Row 3
$arr[1,0] = "ECs"
$arr[1,1] = "Sema3a"
$arr[1,2] = "Plxna2"
$arr[1,3] = "FAPs"
$arr[1,4] = [double]3
$arr[1,5] = [double]1
$arr[1,6] = [double]1.339548666666667
$arr[1,7] = [double]4.018646
$arr[1,8] = [double]0.09827984122213275
$arr[1,9] = [double]0.09827984122213274
$arr[1,10] = [double]3
$arr[1,11] = [double]1
$arr[1,12] = [double]2.577566
$arr[1,13] = [double]7.732698000000001
$arr[1,14] = [double]0.03275155884322009
$arr[1,15] = [double]0.03275155884322008
$arr[1,16] = [double]3.452775098545334
$arr[1,17] = [double]31.07497588690801
$arr[1,18] = [double]0.003218818002889009
$arr[1,19] = [double]0.003218818002889007

# Row 4
$arr[2,0] = "ECs"
$arr[2,1] = "Sema3a"
$arr[2,2] = "Plxna2"
$arr[2,3] = "MuSCs"
$arr[2,4] = [double]3
$arr[2,5] = [double]1
$arr[2,6] = [double]1.339548666666667
$arr[2,7] = [double]4.018646
$arr[2,8] = [double]0.09827984122213275
$arr[2,9] = [double]0.09827984122213274
$arr[2,10] = [double]3
$arr[2,11] = [double]1
$arr[2,12] = [double]35.04673133333333
$arr[2,13] = [double]105.140194
$arr[2,14] = [double]0.4453174364986936
$arr[2,15] = [double]0.4453174364986936
$arr[2,16] = [double]46.94680222859156
$arr[2,17] = [double]422.5212200573241
$arr[2,18] = [double]0.04376572695253879
$arr[2,19] = [double]0.04376572695253878

# Row 5
$arr[3,0] = "ECs"
$arr[3,1] = "Sema3a"
$arr[3,2] = "Plxna2"
$arr[3,3] = "Resolving-Mac"
$arr[3,4] = [double]3
$arr[3,5] = [double]1
$arr[3,6] = [double]1.339548666666667
$arr[3,7] = [double]4.018646
$arr[3,8] = [double]0.09827984122213275
$arr[3,9] = [double]0.09827984122213274
$arr[3,10] = [double]3
$arr[3,11] = [double]1
$arr[3,12] = [double]2.050074333333333
$arr[3,13] = [double]6.150223
$arr[3,14] = [double]0.02604904400552376
$arr[3,15] = [double]0.02604904400552376
$arr[3,16] = [double]2.746174339784222
$arr[3,17] = [double]24.715569058058
$arr[3,18] = [double]0.002560095908851224
$arr[3,19] = [double]0.002560095908851223

# Row 6
$arr[4,0] = "FAPs"
$arr[4,1] = "Sema3a"
$arr[4,2] = "Plxna2"
$arr[4,3] = "ECs"
$arr[4,4] = [double]3
$arr[4,5] = [double]1
$arr[4,6] = [double]2.040291
$arr[4,7] = [double]6.120873
$arr[4,8] = [double]0.1496918182345096
$arr[4,9] = [double]0.1496918182345096
$arr[4,10] = [double]3
$arr[4,11] = [double]1
$arr[4,12] = [double]39.02618766666667
$arr[4,13] = [double]117.078563
$arr[4,14] = [double]0.4958819606525626
$arr[4,15] = [double]0.4958819606525626
$arr[4,16] = [double]79.62477946061099
$arr[4,17] = [double]716.6230151454989
$arr[4,18] = [double]0.07422947231977566
$arr[4,19] = [double]0.07422947231977566

# Row 7
$arr[5,0] = "FAPs"
$arr[5,1] = "Sema3a"
$arr[5,2] = "Plxna2"
$arr[5,3] = "FAPs"
$arr[5,4] = [double]3
$arr[5,5] = [double]1
$arr[5,6] = [double]2.040291
$arr[5,7] = [double]6.120873
$arr[5,8] = [double]0.1496918182345096
$arr[5,9] = [double]0.1496918182345096
$arr[5,10] = [double]3
$arr[5,11] = [double]1
$arr[5,12] = [double]2.577566
$arr[5,13] = [double]7.732698000000001
$arr[5,14] = [double]0.03275155884322009
$arr[5,15] = [double]0.03275155884322008
$arr[5,16] = [double]5.258984711706001
$arr[5,17] = [double]47.330862405354
$arr[5,18] = [double]0.004902640393256149
$arr[5,19] = [double]0.004902640393256148

# Row 8
$arr[6,0] = "FAPs"
$arr[6,1] = "Sema3a"
$arr[6,2] = "Plxna2"
$arr[6,3] = "MuSCs"
$arr[6,4] = [double]3
$arr[6,5] = [double]1
$arr[6,6] = [double]2.040291
$arr[6,7] = [double]6.120873
$arr[6,8] = [double]0.1496918182345096
$arr[6,9] = [double]0.1496918182345096
$arr[6,10] = [double]3
$arr[6,11] = [double]1
$arr[6,12] = [double]35.04673133333333
$arr[6,13] = [double]105.140194
$arr[6,14] = [double]0.4453174364986936
$arr[6,15] = [double]0.4453174364986936
$arr[6,16] = [double]71.50553051881799
$arr[6,17] = [double]643.549774669362
$arr[6,18] = [double]0.06666037676102023
$arr[6,19] = [double]0.06666037676102023

# Row 9
$arr[7,0] = "FAPs"
$arr[7,1] = "Sema3a"
$arr[7,2] = "Plxna2"
$arr[7,3] = "Resolving-Mac"
$arr[7,4] = [double]3
$arr[7,5] = [double]1
$arr[7,6] = [double]2.040291
$arr[7,7] = [double]6.120873
$arr[7,8] = [double]0.1496918182345096
$arr[7,9] = [double]0.1496918182345096
$arr[7,10] = [double]3
$arr[7,11] = [double]1
$arr[7,12] = [double]2.050074333333333
$arr[7,13] = [double]6.150223
$arr[7,14] = [double]0.02604904400552376
$arr[7,15] = [double]0.02604904400552376
$arr[7,16] = [double]4.182748211631
$arr[7,17] = [double]37.644733904679
$arr[7,18] = [double]0.003899328760457605
$arr[7,19] = [double]0.003899328760457605

# Row 10
$arr[8,0] = "MuSCs"
$arr[8,1] = "Sema3a"
$arr[8,2] = "Plxna2"
$arr[8,3] = "ECs"
$arr[8,4] = [double]3
$arr[8,5] = [double]1
$arr[8,6] = [double]10.25010366666667
$arr[8,7] = [double]30.750311
$arr[8,8] = [double]0.7520283405433575
$arr[8,9] = [double]0.7520283405433575
$arr[8,10] = [double]3
$arr[8,11] = [double]1
$arr[8,12] = [double]39.02618766666667
$arr[8,13] = [double]117.078563
$arr[8,14] = [double]0.4958819606525626
$arr[8,15] = [double]0.4958819606525626
$arr[8,16] = [double]400.0224692981214
$arr[8,17] = [double]3600.202223683093
$arr[8,18] = [double]0.3729172879749332
$arr[8,19] = [double]0.3729172879749332

# Row 11
$arr[9,0] = "MuSCs"
$arr[9,1] = "Sema3a"
$arr[9,2] = "Plxna2"
$arr[9,3] = "FAPs"
$arr[9,4] = [double]3
$arr[9,5] = [double]1
$arr[9,6] = [double]10.25010366666667
$arr[9,7] = [double]30.750311
$arr[9,8] = [double]0.7520283405433575
$arr[9,9] = [double]0.7520283405433575
$arr[9,10] = [double]3
$arr[9,11] = [double]1
$arr[9,12] = [double]2.577566
$arr[9,13] = [double]7.732698000000001
$arr[9,14] = [double]0.03275155884322009
$arr[9,15] = [double]0.03275155884322008
$arr[9,16] = [double]26.42031870767534
$arr[9,17] = [double]237.782868369078
$arr[9,18] = [double]0.02463010044707493
$arr[9,19] = [double]0.02463010044707493

# Row 12
$arr[10,0] = "MuSCs"
$arr[10,1] = "Sema3a"
$arr[10,2] = "Plxna2"
$arr[10,3] = "MuSCs"
$arr[10,4] = [double]3
$arr[10,5] = [double]1
$arr[10,6] = [double]10.25010366666667
$arr[10,7] = [double]30.750311
$arr[10,8] = [double]0.7520283405433575
$arr[10,9] = [double]0.7520283405433575
$arr[10,10] = [double]3
$arr[10,11] = [double]1
$arr[10,12] = [double]35.04673133333333
$arr[10,13] = [double]105.140194
$arr[10,14] = [double]0.4453174364986936
$arr[10,15] = [double]0.4453174364986936
$arr[10,16] = [double]359.2326293444816
$arr[10,17] = [double]3233.093664100334
$arr[10,18] = [double]0.3348913327851346
$arr[10,19] = [double]0.3348913327851346

# Row 13
$arr[11,0] = "MuSCs"
$arr[11,1] = "Sema3a"
$arr[11,2] = "Plxna2"
$arr[11,3] = "Resolving-Mac"
$arr[11,4] = [double]3
$arr[11,5] = [double]1
$arr[11,6] = [double]10.25010366666667
$arr[11,7] = [double]30.750311
$arr[11,8] = [double]0.7520283405433575
$arr[11,9] = [double]0.7520283405433575
$arr[11,10] = [double]3
$arr[11,11] = [double]1
$arr[11,12] = [double]2.050074333333333
$arr[11,13] = [double]6.150223
$arr[11,14] = [double]0.02604904400552376
$arr[11,15] = [double]0.02604904400552376
$arr[11,16] = [double]21.01347444103922
$arr[11,17] = [double]189.121269969353
$arr[11,18] = [double]0.01958961933621493
$arr[11,19] = [double]0.01958961933621493

$ws.Range("A2:T13").Value = $arr

# Remove the now-obsolete rows that held the "Resolving-Mac" sending-cluster block
$ws.Range("A14:T17").Delete() | Out-Null

Write-Output "Edit complete"